$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.789941881127277
$ws.Range("C2").Value = 0.1974320777327705
$ws.Range("D2").Value = 0.1255301292643267
$ws.Range("E2").Value = 0.1219734253763445
$ws.Range("F2").Value = 1.673981462271875
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1521136849593776
$ws.Range("L2").Value = 0.3248894851318624
$ws.Range("O2").Value = 4.323744603168876

# Row 3
$ws.Range("B3").Value = 1.658808757678344
$ws.Range("C3").Value = 0.1815793962261125
$ws.Range("D3").Value = 0.1245654263692373
$ws.Range("E3").Value = 0.1228173204742218
$ws.Range("F3").Value = 1.686348976963004
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.153905622380357
$ws.Range("L3").Value = 0.3155631487320676
$ws.Range("O3").Value = 4.363243441035081

# Row 4
$ws.Range("B4").Value = 1.5784718402669
$ws.Range("C4").Value = 0.1718055225642559
$ws.Range("D4").Value = 0.12400262839801
$ws.Range("E4").Value = 0.1233783293661244
$ws.Range("F4").Value = 1.694988679609118
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1550683922273031
$ws.Range("L4").Value = 0.3099255058501171
$ws.Range("O4").Value = 4.390416221686934

# Row 5
$ws.Range("B5").Value = 1.54578101507559
$ws.Range("C5").Value = 0.1678127304095085
$ws.Range("D5").Value = 0.1237807504400408
$ws.Range("E5").Value = 0.1236177360320925
$ws.Range("F5").Value = 1.698772299926425
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1555579589345322
$ws.Range("L5").Value = 0.3076506209181247
$ws.Range("O5").Value = 4.402222931448222

# Row 6
$ws.Range("B6").Value = 1.540355634410162
$ws.Range("C6").Value = 0.1671491427136971
$ws.Range("D6").Value = 0.1237443599413339
$ws.Range("E6").Value = 0.1236581415175309
$ws.Range("F6").Value = 1.699416439484274
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1556402013489588
$ws.Range("L6").Value = 0.3072742422884858
$ws.Range("O6").Value = 4.40422771079281

# Row 7
$ws.Range("B7").Value = 1.578030766879351
$ws.Range("C7").Value = 0.1717517139391873
$ws.Range("D7").Value = 0.1239996057872403
$ws.Range("E7").Value = 0.1233815143758212
$ws.Range("F7").Value = 1.69503864268745
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.155074930998798
$ws.Range("L7").Value = 0.309894734613394
$ws.Range("O7").Value = 4.390572481884249

# Row 8
$ws.Range("B8").Value = 1.744691120654011
$ws.Range("C8").Value = 0.1919745830756199
$ws.Range("D8").Value = 0.1251913936861229
$ws.Range("E8").Value = 0.1222555151122862
$ws.Range("F8").Value = 1.678028590717332
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1527185760664613
$ws.Range("L8").Value = 0.3216554449721087
$ws.Range("O8").Value = 4.336757227981252

# Row 9
$ws.Range("B9").Value = 2.072863257250276
$ws.Range("C9").Value = 0.231302444153954
$ws.Range("D9").Value = 0.1277612453176502
$ws.Range("E9").Value = 0.1203867979884361
$ws.Range("F9").Value = 1.652979266673341
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.1485933163016977
$ws.Range("L9").Value = 0.3454160042077206
$ws.Range("O9").Value = 4.254433096844707

# Row 10
$ws.Range("B10").Value = 2.3147251309681
$ws.Range("C10").Value = 0.2599859456014997
$ws.Range("D10").Value = 0.1297893827025831
$ws.Range("E10").Value = 0.1192198434256806
$ws.Range("F10").Value = 1.639651993101168
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1458638597768207
$ws.Range("L10").Value = 0.3632917006575695
$ws.Range("O10").Value = 4.208149352101543

# Row 11
$ws.Range("B11").Value = 2.424905070806346
$ws.Range("C11").Value = 0.2729871722933979
$ws.Range("D11").Value = 0.130742084007089
$ws.Range("E11").Value = 0.1187335069483701
$ws.Range("F11").Value = 1.63469407715516
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.1446874436800254
$ws.Range("L11").Value = 0.3715135116242863
$ws.Range("O11").Value = 4.19018838638442

# Row 12
$ws.Range("B12").Value = 2.466648134276284
$ws.Range("C12").Value = 0.2779034177767414
$ws.Range("D12").Value = 0.1311071385241576
$ws.Range("E12").Value = 0.1185557308119876
$ws.Range("F12").Value = 1.632975706415777
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.1442513359872928
$ws.Range("L12").Value = 0.3746397025579711
$ws.Range("O12").Value = 4.183832819342882

# Row 13
$ws.Range("B13").Value = 2.457657148297869
$ws.Range("C13").Value = 0.2768449336125798
$ws.Range("D13").Value = 0.1310283274594966
$ws.Range("E13").Value = 0.1185937341337748
$ws.Range("F13").Value = 1.633338709386649
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.1443448426650997
$ws.Range("L13").Value = 0.3739658563562074
$ws.Range("O13").Value = 4.185181758982168

# Row 14
$ws.Range("B14").Value = 2.428338902558892
$ws.Range("C14").Value = 0.2733917773424821
$ws.Range("D14").Value = 0.1307720315677727
$ws.Range("E14").Value = 0.1187187532266876
$ws.Range("F14").Value = 1.634549516179831
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.1446513769610593
$ws.Range("L14").Value = 0.3717704503488335
$ws.Range("O14").Value = 4.189656568711911

# Row 15
$ws.Range("B15").Value = 2.4103832141152
$ws.Range("C15").Value = 0.2712756949907202
$ws.Range("D15").Value = 0.1306156000750107
$ws.Range("E15").Value = 0.1187961626868059
$ws.Range("F15").Value = 1.635311894483408
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.1448403589230969
$ws.Range("L15").Value = 0.3704273583799704
$ws.Range("O15").Value = 4.192455612306475

# Row 16
$ws.Range("B16").Value = 2.307527568170542
$ws.Range("C16").Value = 0.2591353148926032
$ws.Range("D16").Value = 0.1297277235305003
$ws.Range("E16").Value = 0.1192525213793587
$ws.Range("F16").Value = 1.639998253843132
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.145942053525677
$ws.Range("L16").Value = 0.3627561827602221
$ws.Range("O16").Value = 4.209385481660831

# Row 17
$ws.Range("B17").Value = 2.244467362745127
$ws.Range("C17").Value = 0.2516753413961226
$ws.Range("D17").Value = 0.1291907186768739
$ws.Range("E17").Value = 0.119543875575129
$ws.Range("F17").Value = 1.643156291771732
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1466346108370828
$ws.Range("L17").Value = 0.3580731033043776
$ws.Range("O17").Value = 4.220564490935033

# Row 18
$ws.Range("B18").Value = 2.208211584618141
$ws.Range("C18").Value = 0.247380150174223
$ws.Range("D18").Value = 0.1288846824022229
$ws.Range("E18").Value = 0.1197156457013051
$ws.Range("F18").Value = 1.64507668038506
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1470390912267776
$ws.Range("L18").Value = 0.355388008171559
$ws.Range("O18").Value = 4.227285500275428

# Row 19
$ws.Range("B19").Value = 2.195938604262096
$ws.Range("C19").Value = 0.2459251228616779
$ws.Range("D19").Value = 0.1287815518644706
$ws.Range("E19").Value = 0.119774524289145
$ws.Range("F19").Value = 1.645744740810684
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1471770959702914
$ws.Range("L19").Value = 0.354480344541912
$ws.Range("O19").Value = 4.229611091921015

# Row 20
$ws.Range("B20").Value = 2.251178714350374
$ws.Range("C20").Value = 0.252469926896481
$ws.Range("D20").Value = 0.1292475906675818
$ws.Range("E20").Value = 0.1195124267725873
$ws.Range("F20").Value = 1.642809351227342
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.1465602515186157
$ws.Range("L20").Value = 0.3585707479236504
$ws.Range("O20").Value = 4.21934432738027

# Row 21
$ws.Range("B21").Value = 2.436949847284779
$ws.Range("C21").Value = 0.2744062461774206
$ws.Range("D21").Value = 0.1308471958462363
$ws.Range("E21").Value = 0.118681858794524
$ws.Range("F21").Value = 1.634189553571503
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.1445610859891797
$ws.Range("L21").Value = 0.3724149491905138
$ws.Range("O21").Value = 4.188330100251335

# Row 22
$ws.Range("B22").Value = 2.558479460532965
$ws.Range("C22").Value = 0.2887017587523246
$ws.Range("D22").Value = 0.131917601140735
$ws.Range("E22").Value = 0.1181762684098828
$ws.Range("F22").Value = 1.629483387670319
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.1433091618543472
$ws.Range("L22").Value = 0.381537282501327
$ws.Range("O22").Value = 4.170659783744554

# Row 23
$ws.Range("B23").Value = 2.493606756557938
$ws.Range("C23").Value = 0.2810758279868537
$ws.Range("D23").Value = 0.1313440339587828
$ws.Range("E23").Value = 0.1184427087523439
$ws.Range("F23").Value = 1.631910228084507
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.1439723385098626
$ws.Range("L23").Value = 0.3766617784258841
$ws.Range("O23").Value = 4.179852605757702

# Row 24
$ws.Range("B24").Value = 2.248144517197261
$ws.Range("C24").Value = 0.252110714560871
$ws.Range("D24").Value = 0.1292218704438426
$ws.Range("E24").Value = 0.1195266314849963
$ws.Range("F24").Value = 1.642965876643601
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1465938496794594
$ws.Range("L24").Value = 0.3583457401354622
$ws.Range("O24").Value = 4.219895047439024

# Row 25
$ws.Range("B25").Value = 1.983946642053411
$ws.Range("C25").Value = 0.2206995027467826
$ws.Range("D25").Value = 0.1270413041533658
$ws.Range("E25").Value = 0.1208560968123891
$ws.Range("F25").Value = 1.658865163015719
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1496563324308662
$ws.Range("L25").Value = 0.3389141243629723
$ws.Range("O25").Value = 4.274214066712233
